# Edit cv125121a.xlsx: "correção nos dados e inicio da analise PNAD 2009"
#
# The data table (rows 4-26) contained interleaved section-header rows
# (sexo / cor ou raça / grupos de idade / classes de rendimento) that had
# no values of their own, with the real numeric rows below each header
# accidentally holding the PREVIOUS category's numbers (off-by-one data
# entry mistake). This edit removes the 4 section-header rows (and the
# trailing "fonte: ..." row), and re-aligns every remaining label with its
# correct row of figures, compacting the table down to rows 4-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected table: target row number, row label, then values for columns B..H.
$data = @(
    @(4, 'total (1)', 5.387403382729794,11.77063262412606,16.06760978350885,13.30209344255814,14.63828629990292,10.326054480229,12.00031395376143),
    @(5, 'homens', 7.271577638217321,16.53310405266607,21.85371383176248,18.02505052713633,20.62909680823876,12.88354419201064,16.45422638380496),
    @(6, 'mulheres', 7.294352927788966,17.26197342700178,22.56907225174616,19.30354403727936,19.19372751786775,15.04239478523428,15.10216736539094),
    @(7, 'branca', 7.549021017786179,15.90418501959762,21.15296540165263,20.14041167772645,22.73906876758716,15.50021224880463,15.94737500675959),
    @(8, 'preta ou parda', 7.057099679758583,16.11745971058428,23.1114457009257,17.72570656811556,19.05906044136092,13.85817304672361,15.69479699080041),
    @(9, 'preta', 18.09615776116077,59.82834867861214,60.86541468923061,51.69294785629166,32.43510269130071,30.54475345388158,59.14450296003196),
    @(10, 'parda', 7.405512781931758,16.60198252964414,23.65643089906885,18.45729465983731,22.46386332063155,15.24161530995583,16.10745320742356),
    @(11, '10 a 24 anos', 11.94140506999072,23.37734701858705,34.68782407068982,36.62507402402348,31.9203965472073,25.71329318242151,27.75252214305599),
    @(12, '25 a 34 anos', 10.19719656828595,26.13032843535495,37.18774743285265,22.40705268536676,29.11091593137526,18.13987465100105,20.91224131085761),
    @(13, '35 a 39 anos', 13.72647802312186,28.42680176869172,41.4140426950086,42.26528483878885,44.55561055787751,32.08378218406971,33.54907883891536),
    @(14, '40 a 49 anos', 11.8061499979057,30.20394097406351,38.7632851821606,30.34199640306776,30.68105635072406,21.17584209475175,23.36324650148106),
    @(15, '50 anos ou mais', 9.49906606575661,21.81851251253374,28.75690204355052,25.03990364549829,26.86934271084719,19.2473918566034,22.18717737589363),
    @(16, 'total (2)(3)', 5.387403382729794,11.77063262412606,16.06760978350885,13.30209344255814,14.63828629990292,10.326054480229,12.00031395376143),
    @(17, 'sem rendimento a menos de 1/4 do salário mínimo (2) (4) ', 17.859388725917,29.66019362557027,61.73818386071391,49.55528786301879,44.03013636645709,34.32361773428224,40.48719157665269),
    @(18, '1/4 a menos de 1/2 salário mínimo (2)', 12.58311587522816,28.40898233465533,46.47387019501372,37.6371686820284,30.51971707313654,21.98250286363828,31.16152878140587),
    @(19, '1/2 a menos de 1 salário mínimo (2)', 10.39192446351487,25.31944559195118,33.88945735813748,23.92398966870312,27.50105850838492,19.09258228288815,23.40151864763016),
    @(20, '1 a menos de 2 salários mínimos (2) ', 10.19379062859038,21.35129711926809,32.21798808958533,29.10212476557459,31.88352485795141,21.63633018135524,19.06651838522319),
    @(21, '2 salários mínimos ou mais (2) ', 11.7814181367877,40.70473377916795,27.84030921247587,24.08363921615022,32.76276535859011,25.40800766519148,25.27748654287359)
)

foreach ($item in $data) {
    $row = $item[0]
    $label = $item[1]
    $ws.Cells.Item($row, 1).Value = $label
    for ($col = 2; $col -le 8; $col++) {
        $ws.Cells.Item($row, $col).Value = $item[$col]
    }
}

# Drop the now-obsolete trailing rows (old rows 22-26: the last corrected
# data row + the "fonte: ibge..." footer) so the table ends at row 21.
$ws.Range("A22:A26").EntireRow.Delete() | Out-Null
